$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 29.947775
$ws.Range("H2").Value = 89.84332500000001
$ws.Range("I2").Value = 0.9303126840830549
$ws.Range("J2").Value = 0.930312684083055
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 3663.47903204855
$ws.Range("R2").Value = 32971.31128843695
$ws.Range("S2").Value = 0.9101596002940662
$ws.Range("T2").Value = 0.9101596002940664
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 29.947775
$ws.Range("H3").Value = 89.84332500000001
$ws.Range("I3").Value = 0.9303126840830549
$ws.Range("J3").Value = 0.930312684083055
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 11.57054248826667
$ws.Range("R3").Value = 104.1348823944
$ws.Range("S3").Value = 0.002874600955588802
$ws.Range("T3").Value = 0.002874600955588803
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.947775
$ws.Range("H4").Value = 89.84332500000001
$ws.Range("I4").Value = 0.9303126840830549
$ws.Range("J4").Value = 0.930312684083055
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 69.54753819585002
$ws.Range("R4").Value = 625.9278437626501
$ws.Range("S4").Value = 0.01727848283339984
$ws.Range("T4").Value = 0.01727848283339985
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.9818753333333333
$ws.Range("H5").Value = 2.945626
$ws.Range("I5").Value = 0.03050146719708818
$ws.Range("J5").Value = 0.03050146719708818
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 120.1117510650573
$ws.Range("R5").Value = 1081.005759585516
$ws.Range("S5").Value = 0.02984072308962083
$ws.Range("T5").Value = 0.02984072308962084
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9818753333333333
$ws.Range("H6").Value = 2.945626
$ws.Range("I6").Value = 0.03050146719708818
$ws.Range("J6").Value = 0.03050146719708818
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 0.3793547354524445
$ws.Range("R6").Value = 3.414192619072
$ws.Range("S6").Value = 0.00009424739472194757
$ws.Range("T6").Value = 0.00009424739472194757
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.9818753333333333
$ws.Range("H7").Value = 2.945626
$ws.Range("I7").Value = 0.03050146719708818
$ws.Range("J7").Value = 0.03050146719708818
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 2.280203195348
$ws.Range("R7").Value = 20.521828758132
$ws.Range("S7").Value = 0.0005664967127453959
$ws.Range("T7").Value = 0.0005664967127453959
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.261435
$ws.Range("H8").Value = 3.784305
$ws.Range("I8").Value = 0.03918584871985675
$ws.Range("J8").Value = 0.03918584871985676
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.328922
$ws.Range("N8").Value = 366.986766
$ws.Range("O8").Value = 0.9783373008518612
$ws.Range("P8").Value = 0.9783373008518613
$ws.Range("Q8").Value = 154.30998372307
$ws.Range("R8").Value = 1388.78985350763
$ws.Range("S8").Value = 0.03833697746817401
$ws.Range("T8").Value = 0.03833697746817403
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.261435
$ws.Range("H9").Value = 3.784305
$ws.Range("I9").Value = 0.03918584871985675
$ws.Range("J9").Value = 0.03918584871985676
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.3863573333333334
$ws.Range("N9").Value = 1.159072
$ws.Range("O9").Value = 0.003089929874945324
$ws.Range("P9").Value = 0.003089929874945324
$ws.Range("Q9").Value = 0.4873646627733333
$ws.Range("R9").Value = 4.38628196496
$ws.Range("S9").Value = 0.0001210815246345733
$ws.Range("T9").Value = 0.0001210815246345734
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.261435
$ws.Range("H10").Value = 3.784305
$ws.Range("I10").Value = 0.03918584871985675
$ws.Range("J10").Value = 0.03918584871985676
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.322294
$ws.Range("N10").Value = 6.966882000000001
$ws.Range("O10").Value = 0.0185727692731934
$ws.Range("P10").Value = 0.0185727692731934
$ws.Range("Q10").Value = 2.92942293189
$ws.Range("R10").Value = 26.36480638701
$ws.Range("S10").Value = 0.0007277897270481607
$ws.Range("T10").Value = 0.0007277897270481608
